# Generate Report for Handoff
# Adds two new file entries (51eb981c-... and a9e883ab-...) to the
# Overview, zh-cn and de-de sheets, mirroring the existing rows' layout,
# styling (hyperlink font, date number format) and hyperlinks.

$wb = $excel.ActiveWorkbook

$hyperlinkColor = 15570276   # BGR long for ARGB FF6495ED (matches existing "HyperLink" font color)
$dateFormat = "yyyy-mm-dd HH:mm:ss"

function Style-AsHyperlink($range) {
    $range.Font.Underline = 2
    $range.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------
# Sheet "Overview" (columns: File Name | zh-cn | de-de | Latest Handoff Date)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$overviewRows = @(
    @{ Row = 4; File = "51eb981c-dc17-414f-bafa-f65f151ce654.md"; Date = "2016-03-25 07:53:58" },
    @{ Row = 5; File = "a9e883ab-dab4-40f8-b6f9-35275b0ef014.md"; Date = "2016-03-25 07:53:58" }
)

foreach ($r in $overviewRows) {
    $row = $r.Row
    $wsOverview.Range("A$row").Value2 = $r.File
    $wsOverview.Range("B$row").Value2 = "Ready for handoff"
    $wsOverview.Range("C$row").Value2 = "Ready for handoff"
    $wsOverview.Range("D$row").Value2 = $r.Date
    $wsOverview.Range("D$row").NumberFormat = $dateFormat

    $wsOverview.Hyperlinks.Add($wsOverview.Range("A$row"), "https://github.com/OpenLocalizationTest/oltest/blob/02accbff4377db5cbb4d463c8c8bc2770a9ce524/e2e/$($r.File)", "", "", $r.File)
    Style-AsHyperlink $wsOverview.Range("A$row")
}

# ---------------------------------------------------------------------
# Sheet "zh-cn" (columns: A Source File Name | B File Extension | C Status |
#   D Latest Handoff File | E Latest Handoff Datetime | H Latest Handback
#   DateTime | J Handoff Reason)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$zhRows = @(
    @{ Row = 4; File = "51eb981c-dc17-414f-bafa-f65f151ce654.md"; Xlf = "51eb981c-dc17-414f-bafa-f65f151ce654.c6228fe864671987c0d0cc7531906fe94641eda8.zh-cn.xlf"; HandoffDate = "2016-03-25 07:53:49" },
    @{ Row = 5; File = "a9e883ab-dab4-40f8-b6f9-35275b0ef014.md"; Xlf = "a9e883ab-dab4-40f8-b6f9-35275b0ef014.960d03e67802987821fb57b142cd8807c27af897.zh-cn.xlf"; HandoffDate = "2016-03-25 07:53:49" }
)

foreach ($r in $zhRows) {
    $row = $r.Row
    $wsZhCn.Range("A$row").Value2 = $r.File
    $wsZhCn.Range("B$row").Value2 = ".md"
    $wsZhCn.Range("C$row").Value2 = "Ready for handoff"
    $wsZhCn.Range("D$row").Value2 = $r.Xlf
    $wsZhCn.Range("E$row").Value2 = $r.HandoffDate
    $wsZhCn.Range("E$row").NumberFormat = $dateFormat
    $wsZhCn.Range("H$row").Value2 = "0001-01-01 00:00:00"
    $wsZhCn.Range("H$row").NumberFormat = $dateFormat
    $wsZhCn.Range("J$row").Value2 = "Include"

    $wsZhCn.Hyperlinks.Add($wsZhCn.Range("A$row"), "https://github.com/OpenLocalizationTest/oltest/blob/02accbff4377db5cbb4d463c8c8bc2770a9ce524/e2e/$($r.File)", "", "", $r.File)
    Style-AsHyperlink $wsZhCn.Range("A$row")

    $wsZhCn.Hyperlinks.Add($wsZhCn.Range("D$row"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0fa9243e62fde71fb39dbcf252cab93e4dc2f58e/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/$($r.Xlf)", "", "", $r.Xlf)
    Style-AsHyperlink $wsZhCn.Range("D$row")
}

# ---------------------------------------------------------------------
# Sheet "de-de" (same layout as zh-cn, but de-de handoff datetime/xlf files)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$deRows = @(
    @{ Row = 4; File = "51eb981c-dc17-414f-bafa-f65f151ce654.md"; Xlf = "51eb981c-dc17-414f-bafa-f65f151ce654.c6228fe864671987c0d0cc7531906fe94641eda8.de-de.xlf"; HandoffDate = "2016-03-25 07:53:58" },
    @{ Row = 5; File = "a9e883ab-dab4-40f8-b6f9-35275b0ef014.md"; Xlf = "a9e883ab-dab4-40f8-b6f9-35275b0ef014.960d03e67802987821fb57b142cd8807c27af897.de-de.xlf"; HandoffDate = "2016-03-25 07:53:58" }
)

foreach ($r in $deRows) {
    $row = $r.Row
    $wsDeDe.Range("A$row").Value2 = $r.File
    $wsDeDe.Range("B$row").Value2 = ".md"
    $wsDeDe.Range("C$row").Value2 = "Ready for handoff"
    $wsDeDe.Range("D$row").Value2 = $r.Xlf
    $wsDeDe.Range("E$row").Value2 = $r.HandoffDate
    $wsDeDe.Range("E$row").NumberFormat = $dateFormat
    $wsDeDe.Range("H$row").Value2 = "0001-01-01 00:00:00"
    $wsDeDe.Range("H$row").NumberFormat = $dateFormat
    $wsDeDe.Range("J$row").Value2 = "Include"

    $wsDeDe.Hyperlinks.Add($wsDeDe.Range("A$row"), "https://github.com/OpenLocalizationTest/oltest/blob/02accbff4377db5cbb4d463c8c8bc2770a9ce524/e2e/$($r.File)", "", "", $r.File)
    Style-AsHyperlink $wsDeDe.Range("A$row")

    $wsDeDe.Hyperlinks.Add($wsDeDe.Range("D$row"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c1908ca9560353f434faac32bc5d57a98e77b2d2/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/$($r.Xlf)", "", "", $r.Xlf)
    Style-AsHyperlink $wsDeDe.Range("D$row")
}

Write-Host "Added handoff rows to Overview, zh-cn and de-de sheets"
